$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "BackTracking" row from the Easy block (row 5) — it is being
# folded into the renamed "Trees" topic above it (was "Trees and Graphs").
$ws.Rows("5").Delete()

# Rename the Easy-block "Trees and Graphs" topic to "Trees".
$ws.Range("A4").Value = "Trees"

# Add an ETA note next to the Easy-section subtotal.
$ws.Range("F10").Value = "ETA: 10 Days"

# Update selection to match the saved workbook state.
$ws.Range("D3").Select()
